$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 401; this shifts the existing rows
# 401-418 down to 405-422 (Excel copies formatting from the row
# that was previously at 401, same as a normal "Insert" in the UI).
$ws.Rows("401:404").Insert()

# Fill the newly inserted rows with the new price block (date 44706 /
# 2022-05-25), following the same column layout as the rest of the
# sheet: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg,
# Q Kg o Unidades, R Clasificacion.

$newRows = @(
    @{ Row = 401; Calidad = "Especial"; Volumen = 400;  Min = 11000; Max = 12000; Prom = 11500; PrecioKg = 639 },
    @{ Row = 402; Calidad = "Primera";  Volumen = 500;  Min = 9000;  Max = 10000; Prom = 9500;  PrecioKg = 528 },
    @{ Row = 403; Calidad = "Segunda";  Volumen = 360;  Min = 6000;  Max = 7000;  Prom = 6500;  PrecioKg = 361 },
    @{ Row = 404; Calidad = "Tercera";  Volumen = 300;  Min = 4000;  Max = 5000;  Prom = 4500;  PrecioKg = 250 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44706
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = "$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
